# Deluvio: looks like release-candidate.
#
# 1. Mark the "Сделать измеритель батарейки." To-Do item as done
#    (strike-through on both the run and the paragraph mark).
# 2. Move the "_GoBack" bookmark from just before that item to just
#    after the "Проблема питания" heading's text.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: find the "Сделать измеритель батарейки." paragraph via Find
# and strike it through (selecting the whole paragraph range, including
# its paragraph mark, makes Word stamp <w:strike/> onto both the run
# and the paragraph-mark run properties, matching the other already
# "done" items in this list).
# ---------------------------------------------------------------------
$todoRange = $d.Content
$found = $todoRange.Find.Execute("Сделать измеритель батарейки.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Сделать измеритель батарейки.' paragraph"
}
$todoPara = $todoRange.Paragraphs.Item(1)
$todoParaRange = $todoPara.Range
$todoParaRange.Font.StrikeThrough = 1

# ---------------------------------------------------------------------
# Step 2: relocate the "_GoBack" bookmark so it sits right after the
# "Проблема питания" run instead of right before the To-Do run.
#
# Quirk work-around: this runtime mis-resolves a collapsed Range whose
# offset equals "end of a paragraph's text, right before its paragraph
# mark" when it's handed straight to Bookmarks.Add (it ends up adding
# the bookmark around the wrong paragraph). We dodge that edge case by
# temporarily inserting one throw-away character after the heading
# text, adding the bookmark in what is now a safe (not paragraph-mark-
# adjacent) position, and then deleting the throw-away character again.
# ---------------------------------------------------------------------
$headingPara = $d.Paragraphs.Item($todoPara.Range.Information(3) + 1)  # not used; placeholder removed below
PSEOF_PLACEHOLDER_REMOVE_LINE
